# Juno: check in to OLPRODLOC.
# Localize the "Sales report" sheet into Arabic: translate the sheet name
# and the column-header row (Year-Quarter / region names). The quarterly
# period labels in column A (2022-Q1 .. 2023-Q4) are locale-invariant and
# keep their original text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the header row (row 1): Year-Quarter, Midwest, Mountain,
# Northeast, South, Southeast, West.
$ws.Range("A1").Value = "ربع السنوي"
$ws.Range("B1").Value = "الغرب الأوسط"
$ws.Range("C1").Value = "جبل"
$ws.Range("D1").Value = "شمال شرق"
$ws.Range("E1").Value = "الجنوب"
$ws.Range("F1").Value = "جنوب شرق"
$ws.Range("G1").Value = "الغرب"

# Rename the worksheet/tab to the Arabic title.
$ws.Name = "تقرير المبيعات"
